$d = $word.ActiveDocument

# The author justified the body of the interview-conclusions document
# ("Justificación en documento"): every paragraph of the write-up -
# i.e. everything except the three leading "EMPRESA:/PROCEDIMIENTO:/PROCESOS:"
# metadata lines - gets its alignment switched to Justify (wdAlignParagraphJustify = 3).

$wdAlignParagraphJustify = 3
$skipLeadingMetaParagraphs = 3

$count = $d.Paragraphs.Count
for ($i = $skipLeadingMetaParagraphs + 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.ParagraphFormat.Alignment = $wdAlignParagraphJustify
}

Write-Output ("Justified paragraphs " + ($skipLeadingMetaParagraphs + 1) + " through " + $count + " of " + $count)
